$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "a1"
$ws.Range("G3").Value = "3"
$ws.Range("H3").Value = "12"
$ws.Range("I3").Value = "13"
$ws.Range("J3").Value = "df3"
$ws.Range("N3").Value = "b2"
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = "2b"
$ws.Range("R3").Value = "r"

$ws.Range("S12").Select()
